$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("main")

# ---------------------------------------------------------------------------
# 1) Insert 3 blank rows at row 68 (old row 68 "999/test" shifts down to 71)
# ---------------------------------------------------------------------------
$ws.Range("A68:A70").EntireRow.Insert()

# ---------------------------------------------------------------------------
# 2) Free up the shared string currently used by D67 ("BAH Awards") so the
#    pool slot gets reused by the new row's text (matches how the workbook's
#    shared string table ends up ordered in the target file).
# ---------------------------------------------------------------------------
$ws.Cells.Item(67, 4).ClearContents()

# ---------------------------------------------------------------------------
# 3) Row 68 - new "Awards" entry (Booz Allen)
# ---------------------------------------------------------------------------
$ws.Cells.Item(68, 1).Value = 78
$ws.Cells.Item(68, 2).Value = "08 Award"
$ws.Cells.Item(68, 3).Value = 45930
$ws.Cells.Item(68, 4).Value = "73% YoY Reduction in PROD Defects by our Team"
$ws.Cells.Item(68, 5).Value = "[inc]"
$ws.Cells.Item(68, 7).Style = "Normal"
$ws.Cells.Item(68, 9).Style = "Normal"
$ws.Cells.Item(68, 10).Style = "Normal"
$ws.Cells.Item(68, 13).Value = "Booz Allen"
$ws.Cells.Item(68, 16).Value = "Cloud"
$ws.Cells.Item(68, 17).Value = "Cloud"

# ---------------------------------------------------------------------------
# 4) Row 69 - new "Awards" entry (EDHEC)
# ---------------------------------------------------------------------------
$ws.Cells.Item(69, 1).Value = 79
$ws.Cells.Item(69, 2).Value = "08 Award"
$ws.Cells.Item(69, 3).Value = 45930
$ws.Cells.Item(69, 4).Value = "Business School Scholarship Winner"
$ws.Cells.Item(69, 7).Style = "Normal"
$ws.Cells.Item(69, 9).Style = "Normal"
$ws.Cells.Item(69, 10).Style = "Normal"
$ws.Cells.Item(69, 13).Value = "EDHEC"
$ws.Cells.Item(69, 16).Value = "Cloud"
$ws.Cells.Item(69, 17).Value = "Cloud"

# Restyle M69 from the carried-over style (19) to style 22 (vertical top, no wrap)
$ws.Cells.Item(33, 13).Copy()
$ws.Cells.Item(69, 13).PasteSpecial(-4122)

# Build the new "Segoe UI 11 black" font/style once, starting from an existing
# style (fontId 3, no alignment/numberFormat) so only ONE brand-new font is
# created instead of several intermediate throw-away ones.
$ws.Cells.Item(2, 18).Copy()
$ws.Cells.Item(69, 5).PasteSpecial(-4122)
$ws.Cells.Item(69, 5).Value = "[inc]"
$ws.Cells.Item(69, 5).Font.Name = "Segoe UI"

# ---------------------------------------------------------------------------
# 5) Row 70 - new "Awards" entry (NC State)
# ---------------------------------------------------------------------------
$ws.Cells.Item(70, 1).Value = 80
$ws.Cells.Item(70, 2).Value = "08 Award"
$ws.Cells.Item(70, 3).Value = 45930
$ws.Cells.Item(70, 4).Value = "University Activities"
$ws.Cells.Item(70, 7).Style = "Normal"
$ws.Cells.Item(70, 9).Style = "Normal"
$ws.Cells.Item(70, 10).Style = "Normal"
$ws.Cells.Item(70, 13).Value = "NC State"
$ws.Cells.Item(70, 16).Value = "Cloud"
$ws.Cells.Item(70, 17).Value = "Cloud"

$ws.Cells.Item(33, 13).Copy()
$ws.Cells.Item(70, 13).PasteSpecial(-4122)

$ws.Cells.Item(69, 5).Copy()
$ws.Cells.Item(70, 5).PasteSpecial(-4122)
$ws.Cells.Item(70, 5).Value = "[inc]"

# ---------------------------------------------------------------------------
# 6) Row 67 - retitle its D cell from the old "BAH Awards" to "Awards"
# ---------------------------------------------------------------------------
$ws.Cells.Item(67, 4).Value = "Awards"

# ---------------------------------------------------------------------------
# 7) Row height changes
# ---------------------------------------------------------------------------
$ws.Rows.Item(2).RowHeight = 114
$ws.Rows.Item(3).RowHeight = 256.5
$ws.Rows.Item(4).RowHeight = 199.5
$ws.Rows.Item(5).RowHeight = 409.5
$ws.Rows.Item(6).RowHeight = 409.5
$ws.Rows.Item(7).RowHeight = 409.5
$ws.Rows.Item(8).RowHeight = 409.5
$ws.Rows.Item(36).RowHeight = 127.5

# ---------------------------------------------------------------------------
# 8) Selection / active cell bookkeeping
# ---------------------------------------------------------------------------
$ws.Cells.Item(71, 16).Select()
